# Update the Status value on the "Metadata" sheet from "draft" to "active".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B6").Value = "active"
